# Commit: "Corregido los modelos y relaciones"
#
# Adds a new "pagable" worksheet (billetes / pagos / billetables helper
# tables used by the simulation) as the last sheet of the workbook and
# leaves it as the active / selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet as the LAST sheet (after Hoja3) and rename it.
# ---------------------------------------------------------------------
$hoja2 = $wb.Worksheets.Item("Hoja2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "pagable"

# ---------------------------------------------------------------------
# 2. Formats first (so merged "header" cells keep their border/alignment
#    on every cell of the merge, matching how Excel itself lays this
#    out), values afterwards.
# ---------------------------------------------------------------------

# Plain thin-border cells (re-using Hoja2's existing border style, so no
# new border definition is created) over the three table bodies.
$hoja2.Range("B5").Copy()
$ws.Range("A5:H16").PasteSpecial(-4122)
$ws.Range("A22:E28").PasteSpecial(-4122)

# Thin-border + centered cells, used for the three merged table titles.
$ws.Range("A4:B4,E4:H4,A21:E21").HorizontalAlignment = -4108

# The "cantidad" column (D22:D28) uses a border style whose alignment was
# touched (centered then reset) without ever getting a value - reproduce
# it the same way.
$ws.Range("D22:D28").HorizontalAlignment = -4108
$ws.Range("D22:D28").HorizontalAlignment = 1   # xlGeneral

# B23 keeps the workbook's default (no) style.
$ws.Range("B23").ClearFormats()

# ---------------------------------------------------------------------
# 3. Table 1 - BILLETES (id / corte)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "BILLETES"
$ws.Range("A4:B4").Merge()

$ws.Range("A5").Value = "id"
$ws.Range("B5").Value = "corte"

$billetes = @(
    @(1, 200),
    @(2, 100),
    @(3, 50),
    @(4, 20),
    @(5, 10),
    @(6, 5),
    @(7, 2),
    @(8, 1),
    @(9, 0.5),
    @(10, 0.2),
    @(11, 0.1)
)
$r = 6
foreach ($row in $billetes) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $r++
}

# ---------------------------------------------------------------------
# 4. Table 2 - pagos (ID / pagable_type / monto / pagocon)
# ---------------------------------------------------------------------
$ws.Range("E4").Value = "pagos"
$ws.Range("E4:H4").Merge()

$ws.Range("E5").Value = "ID"
$ws.Range("F5").Value = "pagable_type"
$ws.Range("G5").Value = "monto"
$ws.Range("H5").Value = "pagocon"

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "App/Models/Inscripcione"
$ws.Range("G6").Value = 250
$ws.Range("H6").Value = 300

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "App/Models/Inscripcione"
$ws.Range("G7").Value = 50
$ws.Range("H7").Value = 50

# Empty (but still bordered) cells E8:H10 were already formatted above.

# ---------------------------------------------------------------------
# 5. Table 3 - billetables (billete_id / billetable_id / billetable_type /
#    cantidad / c/p)
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "billetables"
$ws.Range("A21:E21").Merge()

$ws.Range("A22").Value = "billete_id"
$ws.Range("B22").Value = "billetable_id"
$ws.Range("C22").Value = "billetable_type"
$ws.Range("D22").Value = "cantidad"
$ws.Range("E22").Value = "c/p"

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "App/Pagos"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1

$ws.Range("A24").Value = 3
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "App/Pagos"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1

$ws.Range("A25").Value = 3
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "App/Pagos"
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0

# Rows 26:28 stay empty (formats only, already applied above).

# ---------------------------------------------------------------------
# 6. Selection / active sheet bookkeeping so tabSelected + activeTab end
#    up on the new sheet, matching the authored selection (D26).
# ---------------------------------------------------------------------
$ws.Range("D26").Select()
$ws.Activate()
